$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9220431447029114
$ws.Range("B1").Value = 1.255401253700256
$ws.Range("C1").Value = 2.110778331756592
$ws.Range("D1").Value = 4.580045700073242
$ws.Range("E1").Value = 2.136010646820068
